# Updated procedure. And added stdev and averages to similarity percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert RS-Sim_STD column before MPM-Pass (old col I) ---
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "RS-Sim_STD"
$ws.Range("I2").Formula = "=STDEV(G2:H2)"
$ws.Range("I3:I9").Formula = "=STDEV(G3:H3)"

# --- 2) Insert MPM-Sim_STD column before MX-Pass (now at col Q) ---
$ws.Columns("Q:Q").Insert()
$ws.Range("Q1").Value = "MPM-Sim_STD"
$ws.Range("Q2").Formula = "=STDEV(O2:P2)"
$ws.Range("Q3:Q9").Formula = "=STDEV(O3:P3)"

# --- 3) Insert MX-Sim-STD column before Exp Order (now at col Y) ---
$ws.Columns("Y:Y").Insert()
$ws.Range("Y1").Value = "MX-Sim-STD"
$ws.Range("Y2").Formula = "=STDEV(W2:X2)"
$ws.Range("Y3:Y9").Formula = "=STDEV(W3:X3)"

# --- 4) Add averages across row 10 for many columns ---
$ws.Range("E10").Formula = "=AVERAGE(E2:E9)"
$ws.Range("G10").Formula = "=AVERAGE(G2:G9)"
$ws.Range("H10").Formula = "=AVERAGE(H2:H9)"
$ws.Range("I10").Formula = "=AVERAGE(I2:I9)"
$ws.Range("M10").Formula = "=AVERAGE(M2:M9)"
$ws.Range("O10").Formula = "=AVERAGE(O2:O9)"
$ws.Range("P10").Formula = "=AVERAGE(P2:P9)"
$ws.Range("U10").Formula = "=AVERAGE(U2:U9)"
$ws.Range("W10").Formula = "=AVERAGE(W2:W9)"
$ws.Range("X10").Formula = "=AVERAGE(X2:X9)"
$ws.Range("Y10").Formula = "=AVERAGE(Y2:Y9)"

# --- 5) Window / view adjustments ---
$excel.ActiveWindow.WindowState = -4143
$ws.Range("F11").Select()
